# Weekly refresh of the Fruta/Hortaliza "Cebollín" dataset:
# each row's Fecha (and, for rows that also moved volume bracket, the
# Volumen/Precio promedio ponderado/Precio $/Kg figures that travel with
# it) is updated to reflect the new reporting week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44321
$ws.Range("J2").Value = 100

# Row 3
$ws.Range("D3").Value = 44321
$ws.Range("J3").Value = 50

# Row 4
$ws.Range("D4").Value = 44525
$ws.Range("O4").Value = "Región de Ñuble"

# Row 5
$ws.Range("D5").Value = 44525
$ws.Range("O5").Value = "Región de Ñuble"

# Row 6
$ws.Range("D6").Value = 44491
$ws.Range("O6").Value = "Región Metropolitana"

# Row 7
$ws.Range("D7").Value = 44491
$ws.Range("O7").Value = "Región Metropolitana"

# Row 8
$ws.Range("D8").Value = 44293

# Row 9
$ws.Range("D9").Value = 44293

# Row 12
$ws.Range("D12").Value = 44328

# Row 13
$ws.Range("D13").Value = 44328

# Row 14
$ws.Range("D14").Value = 44230
$ws.Range("J14").Value = 100

# Row 15
$ws.Range("D15").Value = 44230
$ws.Range("J15").Value = 50

# Row 16
$ws.Range("D16").Value = 44358
$ws.Range("J16").Value = 200

# Row 17
$ws.Range("D17").Value = 44358
$ws.Range("J17").Value = 100

# Row 18
$ws.Range("D18").Value = 44188
$ws.Range("J18").Value = 200
$ws.Range("M18").Value = 650
$ws.Range("P18").Value = 108

# Row 19
$ws.Range("D19").Value = 44188
$ws.Range("J19").Value = 100

# Row 20
$ws.Range("D20").Value = 44554
$ws.Range("J20").Value = 200

# Row 21
$ws.Range("D21").Value = 44554
$ws.Range("J21").Value = 100

# Row 22
$ws.Range("D22").Value = 44335
$ws.Range("J22").Value = 150
$ws.Range("M22").Value = 633
$ws.Range("P22").Value = 106

# Row 23
$ws.Range("D23").Value = 44335
$ws.Range("J23").Value = 50
